# Append the newly digitized "colville" stand-structure rows (qmd, basal_area,
# density) pulled from WebPlotDigitizer into the "data" sheet, rows 1056-1070.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$newRows = @(
    @(45.754, "qmd", "cm", 1, "thinburn", "radcliffe2", "both", "washington", "NA", "commercial", "NA", "NA"),
    @(26.837, "qmd", "cm", 5, "thinburn", "radcliffe2", "both", "washington", "NA", "commercial", "NA", "NA"),
    @(35.635, "qmd", "cm", 10, "thinburn", "radcliffe2", "both", "washington", "NA", "commercial", "NA", "NA"),
    @(36.52, "qmd", "cm", 15, "thinburn", "radcliffe2", "both", "washington", "NA", "commercial", "NA", "NA"),
    @(27.59, "qmd", "cm", 20, "thinburn", "radcliffe2", "both", "washington", "NA", "commercial", "NA", "NA"),
    @(20.211, "basal_area", "m2/ha", 1, "thinburn", "radcliffe2", "both", "washington", "NA", "commercial", "NA", "NA"),
    @(15.531, "basal_area", "m2/ha", 5, "thinburn", "radcliffe2", "both", "washington", "NA", "commercial", "NA", "NA"),
    @(17.89, "basal_area", "m2/ha", 10, "thinburn", "radcliffe2", "both", "washington", "NA", "commercial", "NA", "NA"),
    @(18.229, "basal_area", "m2/ha", 15, "thinburn", "radcliffe2", "both", "washington", "NA", "commercial", "NA", "NA"),
    @(15.946, "basal_area", "m2/ha", 20, "thinburn", "radcliffe2", "both", "washington", "NA", "commercial", "NA", "NA"),
    @(119.383, "density", "trees/ha", 1, "thinburn", "radcliffe2", "both", "washington", "NA", "commercial", "NA", "NA"),
    @(109.326, "density", "trees/ha", 5, "thinburn", "radcliffe2", "both", "washington", "NA", "commercial", "NA", "NA"),
    @(120.034, "density", "trees/ha", 10, "thinburn", "radcliffe2", "both", "washington", "NA", "commercial", "NA", "NA"),
    @(125.571, "density", "trees/ha", 15, "thinburn", "radcliffe2", "both", "washington", "NA", "commercial", "NA", "NA"),
    @(113.831, "density", "trees/ha", 20, "thinburn", "radcliffe2", "both", "washington", "NA", "commercial", "NA", "NA"),
)

$startRow = 1056
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $newRows[$i]
    for ($col = 1; $col -le $rowVals.Count; $col++) {
        $ws.Cells.Item($r, $col).Value = $rowVals[$col - 1]
    }
}

# Match the author's final cursor position (selection) on the sheet.
$ws.Range("H1050").Select()
